$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- Add new time-tracking entries (rows 16-18) ---
# Row 16: 2024-07-11, 09:00 -> 13:20, "Analyse et état de l'art"
$ws.Range("A16").Value = 45484
$ws.Range("B16").Value = 0.375
$ws.Range("C16").Value = 0.55555555555555558
$ws.Range("E16").Value = "Analyse et état de l'art"

# Row 17: 2024-07-12, 10:00 -> 19:25, "Analyse et état de l'art"
$ws.Range("A17").Value = 45485
$ws.Range("B17").Value = 0.41666666666666669
$ws.Range("C17").Value = 0.80902777777777779
$ws.Range("E17").Value = "Analyse et état de l'art"

# Row 18: 2024-07-13, 12:00 -> 17:00, "Analyse et état de l'art"
$ws.Range("A18").Value = 45486
$ws.Range("B18").Value = 0.5
$ws.Range("C18").Value = 0.70833333333333337
$ws.Range("E18").Value = "Analyse et état de l'art"

# --- Update the selected cell to reflect where the user left off ---
$ws.Range("H17").Select() | Out-Null
